$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '42.957.15'
Set-TextValue 'E2' '  -1.69%  '
Set-TextValue 'D3' '2.248.77'
Set-TextValue 'E3' '  -1.58%  '
Set-TextValue 'D5' '116.21'
Set-TextValue 'E5' '  +1.52%  '
Set-TextValue 'D6' '302.10'
Set-TextValue 'E6' '  +13.82%  '
Set-TextValue 'D7' '0.633'
Set-TextValue 'E7' '  -1.60%  '
Set-TextValue 'E8' '  -0.01%  '
Set-TextValue 'D9' '0.624'
Set-TextValue 'E9' '  +1.52%  '
Set-TextValue 'D10' '46.51'
Set-TextValue 'E10' '  -1.56%  '
Set-TextValue 'D11' '0.0936'
Set-TextValue 'E11' '  -0.32%  '
Set-TextValue 'D12' '56.97'
Set-TextValue 'E12' '  +4.60%  '
Set-TextValue 'D13' '9.15'
Set-TextValue 'E13' '  -0.34%  '
Set-TextValue 'D15' '15.42'
Set-TextValue 'E15' '  +0.09%  '
Set-TextValue 'D16' '0.898'
Set-TextValue 'E16' '  +2.88%  '
Set-TextValue 'D17' '2.584.05'
Set-TextValue 'E17' '  -1.74%  '
Set-TextValue 'D18' '2.265.90'
Set-TextValue 'D19' '42.850.64'
Set-TextValue 'E19' '  -1.61%  '
Set-TextValue 'D20' '7.71'
Set-TextValue 'E20' '  +12.54%  '
Set-TextValue 'E21' '  -1.74%  '
Set-TextValue 'D22' '74.23'
Set-TextValue 'E22' '  +2.51%  '
Set-TextValue 'D23' '3.64'
Set-TextValue 'E23' '  +26.47%  '
Set-TextValue 'D24' '2.37'
Set-TextValue 'E24' '  -3.24%  '
Set-TextValue 'D25' '233.38'
Set-TextValue 'E25' '  -1.25%  '
Set-TextValue 'D26' '9.44'
Set-TextValue 'E26' '  -0.35%  '
Set-TextValue 'D27' '12.29'
Set-TextValue 'E27' '  +6.15%  '
Set-TextValue 'E28' '  -1.63%  '
Set-TextValue 'D29' '40.32'
Set-TextValue 'E29' '  -3.60%  '
Set-TextValue 'E30' '  -0.68%  '
Set-TextValue 'E31' '  -3.85%  '
Set-TextValue 'D32' '175.82'
Set-TextValue 'E32' '  +1.33%  '
Set-TextValue 'D33' '21.38'
Set-TextValue 'E33' '  -1.47%  '
Set-TextValue 'D34' '0.0912'
Set-TextValue 'E34' '  +0.21%  '
Set-TextValue 'D35' '4.59'
Set-TextValue 'E35' '  +16.46%  '
Set-TextValue 'D36' '5.67'
Set-TextValue 'E36' '  -0.55%  '
Set-TextValue 'D37' '0.129'
Set-TextValue 'E37' '  -0.94%  '
Set-TextValue 'D38' '4.80'
Set-TextValue 'E38' '  +2.54%  '
Set-TextValue 'D39' '0.0376'
Set-TextValue 'E39' '  -2.08%  '
Set-TextValue 'E40' '  +0.58%  '
Set-TextValue 'D41' '2.62'
Set-TextValue 'E41' '  +2.65%  '
Set-TextValue 'B42' 'Algorand'
Set-TextValue 'C42' 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue 'D42' '0.240'
Set-TextValue 'E42' '  +1.66%  '
Set-TextValue 'B43' 'MultiversX'
Set-TextValue 'C43' 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
Set-TextValue 'D43' '72.68'
Set-TextValue 'E43' '  -1.98%  '
Set-TextValue 'B44' 'Celestia'
Set-TextValue 'C44' 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
Set-TextValue 'D44' '13.54'
Set-TextValue 'E44' '  -5.39%  '
Set-TextValue 'E45' '  +0.14%  '
Set-TextValue 'E46' '  -1.36%  '
Set-TextValue 'D47' '5.61'
Set-TextValue 'E47' '  -5.96%  '
Set-TextValue 'D48' '1.38'
Set-TextValue 'E48' '  +7.44%  '
Set-TextValue 'D49' '108.15'
Set-TextValue 'E49' '  +7.66%  '
Set-TextValue 'D50' '8.65'
Set-TextValue 'E50' '  +0.79%  '
Set-TextValue 'B51' 'TheGraph'
Set-TextValue 'C51' 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextValue 'D51' '0.204'
Set-TextValue 'E51' '  +8.85%  '
